$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet (sheet3): add a formatted-but-empty cell below the existing
# data (stray Hyperlink-style cell at D48), and switch the print paper size
# to A4 (9).
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")
$settings.Activate()
$settings.Cells.Item(48, 4).Style = "Hyperlink"
$settings.PageSetup.PaperSize = 9
$settings.Range("C45").Select()
$excel.ActiveWindow.ScrollRow = 39

# ---------------------------------------------------------------------------
# Companies sheet (sheet4): bring the locator rows in line with the rest of
# the workbook ("Xpath" -> "XPath"), tidy a couple of XPath expressions, fill
# in the address/test-data columns, rename the submit-button row to use an
# Id-type locator, and append two new test rows (lastpage btn / alert).
# ---------------------------------------------------------------------------
$companies = $wb.Worksheets.Item("Companies")

$companies.Range("B2").Value = "XPath"
$companies.Range("C2").Value = ".//*[@id='top-bar-menu']/div[1]/ul/li[6]/ul/li/a"

$companies.Range("B3").Value = "XPath"
$companies.Range("C3").Value = ".//*[@id='top-bar-menu']/div[1]/ul/li[6]/ul/li/ul/li[3]/a"

$companies.Range("B4").Value = "XPath"
$companies.Range("C4").Value = ".//a[@href='/Account/AddNewCompany']"

$companies.Range("B5").Value = "XPath"
$companies.Range("D5").Value = "Newsexport"

$companies.Range("B6").Value = "XPath"

$companies.Range("B7").Value = "XPath"

$companies.Range("B8").Value = "XPath"

$companies.Range("B9").Value = "XPath"

$companies.Range("B10").Value = "XPath"
$companies.Range("D10").Value = "Carlton gore rd"

$companies.Range("B11").Value = "XPath"
$companies.Range("D11").Value = "Newmarket"

$companies.Range("B12").Value = "XPath"
$companies.Range("D12").Value = "Auckland"

$companies.Range("B13").Value = "XPath"
$companies.Range("D13").Value = 1023

$companies.Range("B14").Value = "XPath"

$companies.Range("B15").Value = "XPath"

$companies.Range("B16").Value = "Id"
$companies.Range("C16").Value = "submitBtn"

$companies.Range("A17").Value = "lastpage btn"
$companies.Range("B17").Value = "XPath"
$companies.Range("C17").Value = ".//*[@id='btn4']"

$companies.Range("A18").Value = "alert"
$companies.Range("B18").Value = "XPath"
$companies.Range("C18").Value = ".//*[@id='beehive-alert']/p"

$companies.PageSetup.PaperSize = 9

# Leave the workbook with the Companies sheet active/selected, matching the
# final cursor position recorded in the saved file.
$companies.Activate()
$companies.Range("C16").Select()
